$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The emotion-detection categories changed ("angry" dropped, "fear" and
# "surprise" introduced) and the elapsed-time readings were recomputed
# because the app now has a button to run without showing the webcam
# preview of detected faces (faster loop -> smaller elapsed times).

# Existing rows: refresh the detected-label text and the elapsed time
$ws.Range("B3").Value = "fear"
$ws.Range("C3").Value = 1.02808640000876

$ws.Range("B4").Value = "neutral"
$ws.Range("C4").Value = 1.24741820001509

$ws.Range("B5").Value = "surprise"
$ws.Range("C5").Value = 2.994884000014281

$ws.Range("B6").Value = "neutral"
$ws.Range("C6").Value = 3.202956299996004

# New detection rows captured during the run
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "happy"
$ws.Range("C7").Value = 3.436031999997795

$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "neutral"
$ws.Range("C8").Value = 4.098535500001162

$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "happy"
$ws.Range("C9").Value = 4.34527469999739

# Give the new A-column cells the same style (bordered, bold, centered) as
# the rest of the index column
$ws.Range("A2").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$ws.Range("A9").PasteSpecial(-4122)
